# Scheduled price/profit refresh for the Leve profit tracker workbook.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ,
# LevePriceNQ / LevePriceHQ and the derived LeveProfitNQ / LeveProfitHQ columns
# (H, I, J, K, L, M, N) for the rows whose market data changed, across all
# eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 2244.5
$ws.Range("I43").Value = 2244.5
$ws.Range("K43").Value = 2244.5
$ws.Range("M43").Value = -2175.5

# Row 54
$ws.Range("H54").Value = 19999.5
$ws.Range("I54").Value = 19999.5
$ws.Range("K54").Value = 19999.5
$ws.Range("M54").Value = -19513.5

# Row 80
$ws.Range("H80").Value = 562.4286
$ws.Range("I80").Value = 141.66667
$ws.Range("K80").Value = 425.00001
$ws.Range("M80").Value = 572.99999

# Row 83
$ws.Range("H83").Value = 562.4286
$ws.Range("I83").Value = 141.66667
$ws.Range("K83").Value = 1275.00003
$ws.Range("M83").Value = 3716.99997

# Row 93
$ws.Range("H93").Value = 46662.668
$ws.Range("J93").Value = 46662.668
$ws.Range("L93").Value = 46662.668
$ws.Range("N93").Value = -51654.668

# Row 98
$ws.Range("H98").Value = 7500
$ws.Range("I98").Value = 5000
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 5000
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -3502
$ws.Range("N98").Value = -12996

# Row 122
$ws.Range("H122").Value = 7500
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -34900

# Row 132
$ws.Range("H132").Value = 2003.7542
$ws.Range("I132").Value = 2022.4
$ws.Range("K132").Value = 6067.200000000001
$ws.Range("M132").Value = -3537.200000000001

# Row 137
$ws.Range("H137").Value = 1972.2
$ws.Range("I137").Value = 2108.6667
$ws.Range("K137").Value = 6326.000100000001
$ws.Range("M137").Value = -3776.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 609.06665
$ws.Range("I2").Value = 474.7143
$ws.Range("K2").Value = 474.7143
$ws.Range("M2").Value = -361.7143

# Row 44
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976

# Row 61
$ws.Range("H61").Value = 1005254.4
$ws.Range("I61").Value = 3980
$ws.Range("K61").Value = 3980
$ws.Range("M61").Value = -3768

# Row 94
$ws.Range("H94").Value = 22666.666
$ws.Range("J94").Value = 22666.666
$ws.Range("L94").Value = 22666.666
$ws.Range("N94").Value = -24468.666

# Row 116
$ws.Range("H116").Value = 609.06665
$ws.Range("I116").Value = 474.7143
$ws.Range("K116").Value = 474.7143
$ws.Range("M116").Value = 1819.2857

# Row 132
$ws.Range("H132").Value = 3709.7856
$ws.Range("I132").Value = 3456.6924
$ws.Range("K132").Value = 10370.0772
$ws.Range("M132").Value = -7840.0772

# Row 133
$ws.Range("H133").Value = 134494.67
$ws.Range("J133").Value = 134494.67
$ws.Range("L133").Value = 134494.67
$ws.Range("N133").Value = -139554.67

# Row 136
$ws.Range("H136").Value = 1005254.4
$ws.Range("I136").Value = 3980
$ws.Range("K136").Value = 11940
$ws.Range("M136").Value = -9390

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 609.06665
$ws.Range("I3").Value = 474.7143
$ws.Range("K3").Value = 474.7143
$ws.Range("M3").Value = -360.7143

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3201.6
$ws.Range("I58").Value = 2237.4443
$ws.Range("K58").Value = 2237.4443
$ws.Range("M58").Value = -2034.4443

# Row 94
$ws.Range("H94").Value = 2643.2
$ws.Range("I94").Value = 2208.25
$ws.Range("K94").Value = 2208.25
$ws.Range("M94").Value = -1757.25

# Row 122
$ws.Range("H122").Value = 1199.3684
$ws.Range("J122").Value = 1300
$ws.Range("L122").Value = 3900
$ws.Range("N122").Value = -8800

# Row 134
$ws.Range("H134").Value = 8499.594999999999
$ws.Range("I134").Value = 4921.8335
$ws.Range("J134").Value = 23832.857
$ws.Range("K134").Value = 14765.5005
$ws.Range("L134").Value = 71498.571
$ws.Range("M134").Value = -12230.5005
$ws.Range("N134").Value = -76568.571

# Row 136
$ws.Range("H136").Value = 3201.6
$ws.Range("I136").Value = 2237.4443
$ws.Range("K136").Value = 6712.3329
$ws.Range("M136").Value = -4162.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4998.857
$ws.Range("I3").Value = 2199.4
$ws.Range("J3").Value = 11997.5
$ws.Range("K3").Value = 6598.200000000001
$ws.Range("L3").Value = 35992.5
$ws.Range("M3").Value = -6486.200000000001
$ws.Range("N3").Value = -36216.5

# Row 81
$ws.Range("H81").Value = 4920.2354
$ws.Range("I81").Value = 2531.8
$ws.Range("K81").Value = 7595.400000000001
$ws.Range("M81").Value = -6472.400000000001

# Row 84
$ws.Range("H84").Value = 4920.2354
$ws.Range("I84").Value = 2531.8
$ws.Range("K84").Value = 22786.2
$ws.Range("M84").Value = -17170.2

# Row 113
$ws.Range("H113").Value = 2632266.8
$ws.Range("I113").Value = 4546200.5
$ws.Range("J113").Value = 608.375
$ws.Range("K113").Value = 13638601.5
$ws.Range("L113").Value = 1825.125
$ws.Range("M113").Value = -13636431.5
$ws.Range("N113").Value = -6165.125

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3163.375
$ws.Range("I80").Value = 3219.8
$ws.Range("J80").Value = 3069.3333
$ws.Range("K80").Value = 3219.8
$ws.Range("L80").Value = 3069.3333
$ws.Range("M80").Value = -2221.8
$ws.Range("N80").Value = -5065.3333

# Row 83
$ws.Range("H83").Value = 3163.375
$ws.Range("I83").Value = 3219.8
$ws.Range("J83").Value = 3069.3333
$ws.Range("K83").Value = 16099
$ws.Range("L83").Value = 15346.6665
$ws.Range("M83").Value = -11107
$ws.Range("N83").Value = -25330.6665

# Row 93
$ws.Range("H93").Value = 29899.5
$ws.Range("J93").Value = 29899.5
$ws.Range("L93").Value = 29899.5
$ws.Range("N93").Value = -33643.5

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

# Row 107
$ws.Range("H107").Value = 431.03226
$ws.Range("I107").Value = 519.5263
$ws.Range("K107").Value = 519.5263
$ws.Range("M107").Value = 1400.4737

# Row 132
$ws.Range("H132").Value = 4016.6667
$ws.Range("I132").Value = 3581.25
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 10743.75
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -8213.75
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 225.78947
$ws.Range("I55").Value = 205
$ws.Range("K55").Value = 205
$ws.Range("M55").Value = -32

# Row 61
$ws.Range("H61").Value = 1940.8387
$ws.Range("I61").Value = 1790.75
$ws.Range("K61").Value = 1790.75
$ws.Range("M61").Value = -1588.75

# Row 68
$ws.Range("H68").Value = 3370.25
$ws.Range("I68").Value = 3222.1
$ws.Range("J68").Value = 4111
$ws.Range("K68").Value = 3222.1
$ws.Range("L68").Value = 4111
$ws.Range("M68").Value = -2473.1
$ws.Range("N68").Value = -5609

# Row 71
$ws.Range("H71").Value = 3370.25
$ws.Range("I71").Value = 3222.1
$ws.Range("J71").Value = 4111
$ws.Range("K71").Value = 16110.5
$ws.Range("L71").Value = 20555
$ws.Range("M71").Value = -12366.5
$ws.Range("N71").Value = -28043

# Row 82
$ws.Range("H82").Value = 7764.263
$ws.Range("I82").Value = 11265.818
$ws.Range("J82").Value = 2949.625
$ws.Range("K82").Value = 11265.818
$ws.Range("L82").Value = 2949.625
$ws.Range("M82").Value = -10904.818
$ws.Range("N82").Value = -3671.625

# Row 85
$ws.Range("H85").Value = 7764.263
$ws.Range("I85").Value = 11265.818
$ws.Range("J85").Value = 2949.625
$ws.Range("K85").Value = 11265.818
$ws.Range("L85").Value = 2949.625
$ws.Range("M85").Value = -10017.818
$ws.Range("N85").Value = -5445.625

# Row 113
$ws.Range("H113").Value = 1940.8387
$ws.Range("I113").Value = 1790.75
$ws.Range("K113").Value = 1790.75
$ws.Range("M113").Value = 379.25

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 90833.164
$ws.Range("J46").Value = 90833.164
$ws.Range("L46").Value = 90833.164
$ws.Range("N46").Value = -91295.164

# Row 54
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51040

# Row 100
$ws.Range("H100").Value = 912.1875
$ws.Range("I100").Value = 907
$ws.Range("K100").Value = 1814
$ws.Range("M100").Value = -1273

# Row 132
$ws.Range("H132").Value = 2615.8286
$ws.Range("I132").Value = 2415.3667
$ws.Range("K132").Value = 7246.1001
$ws.Range("M132").Value = -4716.1001

# Row 134
$ws.Range("H134").Value = 90833.164
$ws.Range("J134").Value = 90833.164
$ws.Range("L134").Value = 272499.492
$ws.Range("N134").Value = -277569.492
